$d = $word.ActiveDocument
$st = $d.Styles("Heading 1")
$st.NameLocal = "Overskrift1"
Write-Output $st.NameLocal
